$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New class codes (column A)
$classCodes = @("ps09125","ps09126","ps09127","ps09128","ps09129","ps09130","ps09131","ps09132","ps09133","ps09134","ps09135","ps09136","ps09137","ps09138","ps09139")
# New names (column B)
$names = @("Đen Vâu","Phương Mỹ Chi","Khoa Phạm","Code X","Điệp Vấn 2","Rô nal đô","Tiki Nguyễn","Ô Lông Nguyễn","Heo Xinh Trần","Hà Anh Tuấn","Sơn Tùng","What đờ phắc","Ta đa nguyễn","Mèo ú ","Mèo mướp")

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $classCodes[$i]
    $ws.Cells.Item($row, 2).Value = $names[$i]
    # Column E keeps its quote-prefix ('text') style, so re-assert the
    # leading apostrophe when writing the new class code value.
    $ws.Cells.Item($row, 5).Value = "'WD14306"
}

# New cell I10 = " " (a blank/dash placeholder)
$ws.Cells.Item(10, 9).Value = " "

# Update sheet view: remove top-left scroll anchor, move selection
$ws.Range("G12").Select()
